$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OTIS")

# Row 4: Inventory
$ws.Range("B4").Value = 686000000.0
$ws.Range("C4").Value = 659000000.0
$ws.Range("D4").Value = 667000000.0
$ws.Range("E4").Value = 629000000.0
$ws.Range("F4").Value = 599000000.0

# Row 15: Accounts Payable
$ws.Range("B15").Value = 1459000000.0
$ws.Range("C15").Value = 1453000000.0
$ws.Range("D15").Value = 1392000000.0
$ws.Range("E15").Value = 1349000000.0
$ws.Range("F15").Value = 1102000000.0

# Row 41: Net Debt
$ws.Range("G41").Value = -1407000000.0

# Row 42: Total Debt
$ws.Range("G42").Value = 39000000.0
